# T620X_IO.xlsx WIP "sonsite small plc hmi fix ups" edits
# -------------------------------------------------------
# All of the substantive content changes live on the SNAGS sheet (sheet6.xml):
# a new INPUTS snag row is inserted, several existing snags are annotated with
# DONE / SKIP / FIXED status, and a batch of new snags (PINCODE, UDP,
# HYD_SETUP, TRACKS, RADIO, FAULTS, ECU, NAV) is appended further down the
# sheet. The SNAGS tab also becomes the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNAGS")

# Insert a blank row at 14 to make room for the new INPUTS entry ("machine
# mimic psi or bar?"). This pushes the existing OUTPUTS/POWER UP/DM1/J1939/
# MACHINE_APP/SETTINGS/MIMICS rows down by one, which is exactly the target
# layout (old row 15 -> new row 16, old row 33 -> new row 34, etc).
$ws.Rows.Item(14).Insert()

# --- Mark existing INPUTS diagnostics snags as DONE ---
$ws.Range("D7").Value = "DONE"
$ws.Range("D8").Value = "DONE"
$ws.Range("D9").Value = "DONE"
$ws.Range("D10").Value = "DONE"
$ws.Range("D11").Value = "DONE"
$ws.Range("D12").Value = "DONE"

# --- New INPUTS snag in the row opened up above ---
$ws.Range("B14").Value = "INPUTS"
$ws.Range("C14").Value = "machine mimic psi or bar?"
$ws.Range("D14").Value = "DONE"
$ws.Range("E14").Value = "psi!"

# --- Mark DM1 snags ---
$ws.Range("D20").Value = "DONE"
$ws.Range("D21").Value = "DONE"
$ws.Range("D22").Value = "SKIP"

# --- New snags appended near the bottom of the sheet ---
$ws.Range("B37").Value = "PINCODE"
$ws.Range("C37").Value = "4 digits message instead of 6"

$ws.Range("B40").Value = "UDP"
$ws.Range("C40").Value = "PLC warning and error timeouts"

$ws.Range("B42").Value = "HYD_SETUP"
$ws.Range("C42").Value = "indication that setup mode is on / starting"

$ws.Range("B44").Value = "TRACKS"
$ws.Range("C44").Value = "radio stop pressed - engine still runs"

$ws.Range("B45").Value = "RADIO"
$ws.Range("C45").Value = "very slow interlock on jacks/side - long delay?"

$ws.Range("B46").Value = "FAULTS"
$ws.Range("C46").Value = "DM1 's not added to the fault logger"

$ws.Range("B48").Value = "ECU"
$ws.Range("C48").Value = "oil pressure not showing"
$ws.Range("D48").Value = "FIXED"

$ws.Range("B50").Value = "NAV"
$ws.Range("C50").Value = "IO first page entry index issue"
$ws.Range("E50").Value = "left / right arrows"

$ws.Range("B51").Value = "NAV"
$ws.Range("C51").Value = "engine first page entry index issue"
$ws.Range("E51").Value = "left / right arrows"

# Column C was widened (and its old best-fit auto width dropped) to fit the
# new, longer snag descriptions.
$ws.Columns.Item(3).ColumnWidth = 66.57142857142857

# Selection moved to the last entered cell.
$ws.Range("D51").Select() | Out-Null

# SNAGS becomes the active/selected sheet tab (this also clears
# tabSelected on whichever sheet was active before, e.g. IO).
$ws.Activate() | Out-Null
